# Natmi following Dr Hou advice
# Updates the LR-pair table (Sema3a-Nrp1) with the recomputed per-cluster-pair
# statistics, adding an "ECs" target-cluster series and extending the table
# from 6 to 9 data rows (3 sending clusters x 3 target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 9,20

$arr[0,0] = "ECs"
$arr[0,1] = "Sema3a"
$arr[0,2] = "Nrp1"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 1.397441
$arr[0,7] = 4.192323
$arr[0,8] = 0.6676161521996591
$arr[0,9] = 0.6676161521996592
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 52.91030366666666
$arr[0,13] = 158.730911
$arr[0,14] = 0.4161415425564564
$arr[0,15] = 0.4161415425564564
$arr[0,16] = 73.93902766625033
$arr[0,17] = 665.451248996253
$arr[0,18] = 0.2778228154119721
$arr[0,19] = 0.2778228154119721

$arr[1,0] = "ECs"
$arr[1,1] = "Sema3a"
$arr[1,2] = "Nrp1"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 1.397441
$arr[1,7] = 4.192323
$arr[1,8] = 0.6676161521996591
$arr[1,9] = 0.6676161521996592
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 47.26005833333333
$arr[1,13] = 141.780175
$arr[1,14] = 0.3717021489810786
$arr[1,15] = 0.3717021489810786
$arr[1,16] = 66.04314317739166
$arr[1,17] = 594.388288596525
$arr[1,18] = 0.2481543584670921
$arr[1,19] = 0.2481543584670921

$arr[2,0] = "ECs"
$arr[2,1] = "Sema3a"
$arr[2,2] = "Nrp1"
$arr[2,3] = "sCs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 1.397441
$arr[2,7] = 4.192323
$arr[2,8] = 0.6676161521996591
$arr[2,9] = 0.6676161521996592
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 26.97460733333333
$arr[2,13] = 80.923822
$arr[2,14] = 0.2121563084624651
$arr[2,15] = 0.2121563084624651
$arr[2,16] = 37.69542224650067
$arr[2,17] = 339.258800218506
$arr[2,18] = 0.1416389783205949
$arr[2,19] = 0.141638978320595

$arr[3,0] = "FAPs"
$arr[3,1] = "Sema3a"
$arr[3,2] = "Nrp1"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 0.2347746666666667
$arr[3,7] = 0.7043240000000001
$arr[3,8] = 0.112161700990566
$arr[3,9] = 0.112161700990566
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 52.91030366666666
$arr[3,13] = 158.730911
$arr[3,14] = 0.4161415425564564
$arr[3,15] = 0.4161415425564564
$arr[3,16] = 12.42199890657378
$arr[3,17] = 111.797990159164
$arr[3,18] = 0.04667514326597016
$arr[3,19] = 0.04667514326597017

$arr[4,0] = "FAPs"
$arr[4,1] = "Sema3a"
$arr[4,2] = "Nrp1"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.2347746666666667
$arr[4,7] = 0.7043240000000001
$arr[4,8] = 0.112161700990566
$arr[4,9] = 0.112161700990566
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 47.26005833333333
$arr[4,13] = 141.780175
$arr[4,14] = 0.3717021489810786
$arr[4,15] = 0.3717021489810786
$arr[4,16] = 11.09546444185555
$arr[4,17] = 99.8591799767
$arr[4,18] = 0.04169074529156656
$arr[4,19] = 0.04169074529156656

$arr[5,0] = "FAPs"
$arr[5,1] = "Sema3a"
$arr[5,2] = "Nrp1"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 0.2347746666666667
$arr[5,7] = 0.7043240000000001
$arr[5,8] = 0.112161700990566
$arr[5,9] = 0.112161700990566
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 26.97460733333333
$arr[5,13] = 80.923822
$arr[5,14] = 0.2121563084624651
$arr[5,15] = 0.2121563084624651
$arr[5,16] = 6.332954445147556
$arr[5,17] = 56.99659000632801
$arr[5,18] = 0.02379581243302931
$arr[5,19] = 0.02379581243302931

$arr[6,0] = "sCs"
$arr[6,1] = "Sema3a"
$arr[6,2] = "Nrp1"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 0.4609646666666667
$arr[6,7] = 1.382894
$arr[6,8] = 0.2202221468097748
$arr[6,9] = 0.2202221468097748
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 52.91030366666666
$arr[6,13] = 158.730911
$arr[6,14] = 0.4161415425564564
$arr[6,15] = 0.4161415425564564
$arr[6,16] = 24.38978049293711
$arr[6,17] = 219.508024436434
$arr[6,18] = 0.09164358387851407
$arr[6,19] = 0.09164358387851408

$arr[7,0] = "sCs"
$arr[7,1] = "Sema3a"
$arr[7,2] = "Nrp1"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 0.4609646666666667
$arr[7,7] = 1.382894
$arr[7,8] = 0.2202221468097748
$arr[7,9] = 0.2202221468097748
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 47.26005833333333
$arr[7,13] = 141.780175
$arr[7,14] = 0.3717021489810786
$arr[7,15] = 0.3717021489810786
$arr[7,16] = 21.78521703627222
$arr[7,17] = 196.06695332645
$arr[7,18] = 0.08185704522241986
$arr[7,19] = 0.08185704522241988

$arr[8,0] = "sCs"
$arr[8,1] = "Sema3a"
$arr[8,2] = "Nrp1"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 0.4609646666666667
$arr[8,7] = 1.382894
$arr[8,8] = 0.2202221468097748
$arr[8,9] = 0.2202221468097748
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 26.97460733333333
$arr[8,13] = 80.923822
$arr[8,14] = 0.2121563084624651
$arr[8,15] = 0.2121563084624651
$arr[8,16] = 12.43434087787422
$arr[8,17] = 111.909067900868
$arr[8,18] = 0.04672151770884086
$arr[8,19] = 0.04672151770884087

$ws.Range("A2:T10").Value2 = $arr

